$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate header row to Russian
$ws.Range("A1").Value = "Название"
$ws.Range("B1").Value = "Норма азота"
$ws.Range("C1").Value = "Норма фосфора"
$ws.Range("D1").Value = "Норма калия"
$ws.Range("E1").Value = "Id культуры"
$ws.Range("F1").Value = "Район"
$ws.Range("G1").Value = "Цена"
$ws.Range("H1").Value = "Описание "
$ws.Range("I1").Value = "Назначение"

# Update data values
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 123
$ws.Range("C3").Value = 1

# Column widths (target widths quantized to the nearest value this
# engine's pixel-grid ColumnWidth rounding can actually reach)
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(8).ColumnWidth = 15.0
$ws.Columns.Item(9).ColumnWidth = 16.833333333333332

# Selection change
$ws.Range("J1").Select()
